$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# "Import Priorities" sheet: insert a new row for the Strategic Reserve
# Operator entry (previously its own sheet) above "TechnologiesEmlab".
$ws1 = $wb.Worksheets.Item("Import Priorities")
$ws1.Rows.Item(6).Insert() | Out-Null
$ws1.Cells.Item(6, 1).Value = "StrategicReserveOperator"
$ws1.Cells.Item(6, 2).Value = 7

# Remove the now-redundant "StrategicReserveOperators" sheet.
$ws3 = $wb.Worksheets.Item("StrategicReserveOperators")
$ws3.Delete() | Out-Null

# Make "Import Priorities" the active sheet/selection.
$ws1.Activate() | Out-Null
$ws1.Range("F8").Select() | Out-Null
